$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - copy style from G1 ("sum") and set value to "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
